$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run averaged-intensity computations now that three additional orientation
# schemes (Spiral-90deg-10rot-5space, Spiral-90deg-15rot-5space,
# Spiral-90deg-10rot-3space) were added to the study. Gaussian-Quadrature moves up
# next to the other "whole pole figure" schemes, and the existing rotation/grid
# schemes shift down to make room, plus the hex-grid rows reappear at the bottom
# with their (lightly recomputed) values.

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "Gaussian-Quadrature"
$ws.Range("C10").Value = 0.995896977005419
$ws.Range("D10").Value = 0.9884528854690178
$ws.Range("E10").Value = 0.9946551405311592
$ws.Range("F10").Value = 0.995896977005419
$ws.Range("G10").Value = 0.9812925782960419
$ws.Range("H10").Value = 0.9923753196759012
$ws.Range("I10").Value = 0.9906682252443847
$ws.Range("J10").Value = 0.9884528854690178
$ws.Range("K10").Value = 0.9915540130000885
$ws.Range("L10").Value = 0.9937254950027536
$ws.Range("M10").Value = 0.9905568543703205

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "Spiral-90deg-10rot-5space"
$ws.Range("C11").Value = 0.9916944355013595
$ws.Range("D11").Value = 0.9771689459665518
$ws.Range("E11").Value = 0.9960292800793002
$ws.Range("F11").Value = 0.9916944355013595
$ws.Range("G11").Value = 0.9825779207749094
$ws.Range("H11").Value = 1.007243282868975
$ws.Range("I11").Value = 0.99568469706686
$ws.Range("J11").Value = 0.9771689459665518
$ws.Range("K11").Value = 0.986599113022926
$ws.Range("L11").Value = 0.9891467742621427
$ws.Range("M11").Value = 0.9917330937096592

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "Spiral-90deg-15rot-5space"
$ws.Range("C12").Value = 0.9916599013470611
$ws.Range("D12").Value = 0.977385247141883
$ws.Range("E12").Value = 0.9959866672158417
$ws.Range("F12").Value = 0.9916599013470611
$ws.Range("G12").Value = 0.9828133310953702
$ws.Range("H12").Value = 1.007036198226405
$ws.Range("I12").Value = 0.99568617313607
$ws.Range("J12").Value = 0.977385247141883
$ws.Range("K12").Value = 0.9866859571788624
$ws.Range("L12").Value = 0.9891729292629619
$ws.Range("M12").Value = 0.991761253027105

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "Spiral-90deg-10rot-3space"
$ws.Range("C13").Value = 0.9916753340548542
$ws.Range("D13").Value = 0.9771968529857997
$ws.Range("E13").Value = 0.99598943056112
$ws.Range("F13").Value = 0.9916753340548542
$ws.Range("G13").Value = 0.9826423952537202
$ws.Range("H13").Value = 1.007239599176238
$ws.Range("I13").Value = 0.9956982663745929
$ws.Range("J13").Value = 0.9771968529857997
$ws.Range("K13").Value = 0.9865931417734599
$ws.Range("L13").Value = 0.989134237914157
$ws.Range("M13").Value = 0.9917403130677207

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NoRotation-tilt60deg"
$ws.Range("C14").Value = 0.9991759999999993
$ws.Range("D14").Value = 0.9760199999999979
$ws.Range("E14").Value = 0.9991759999999991
$ws.Range("F14").Value = 0.9991759999999993
$ws.Range("G14").Value = 0.9797720000000018
$ws.Range("H14").Value = 0.9918840000000002
$ws.Range("I14").Value = 0.9947960000000008
$ws.Range("J14").Value = 0.9760199999999979
$ws.Range("K14").Value = 0.9875979999999984
$ws.Range("L14").Value = 0.9933869999999989
$ws.Range("M14").Value = 0.9901373333333332

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "Rotation-NoTilt"
$ws.Range("C15").Value = 1.01
$ws.Range("D15").Value = 0.9438874999999985
$ws.Range("E15").Value = 1.01
$ws.Range("F15").Value = 1.01
$ws.Range("G15").Value = 0.96
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = 1
$ws.Range("J15").Value = 0.9438874999999985
$ws.Range("K15").Value = 0.9769437499999993
$ws.Range("L15").Value = 0.9934718749999996
$ws.Range("M15").Value = 0.9873145833333331

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Rotation-60detTilt"
$ws.Range("C16").Value = 1.002721036083202
$ws.Range("D16").Value = 0.963919763251197
$ws.Range("E16").Value = 1.002259147571201
$ws.Range("F16").Value = 1.002721036083202
$ws.Range("G16").Value = 0.9737184090111993
$ws.Range("H16").Value = 0.9974377936895994
$ws.Range("I16").Value = 0.9972639737856008
$ws.Range("J16").Value = 0.963919763251197
$ws.Range("K16").Value = 0.9830894554111989
$ws.Range("L16").Value = 0.9929052457472004
$ws.Range("M16").Value = 0.9895533538986666

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C17").Value = 0.9923393866502159
$ws.Range("D17").Value = 0.9924990944342341
$ws.Range("E17").Value = 0.9922505098929241
$ws.Range("F17").Value = 0.9923393866502159
$ws.Range("G17").Value = 0.9917670330760217
$ws.Range("H17").Value = 0.9924019427447981
$ws.Range("I17").Value = 0.9928189273290116
$ws.Range("J17").Value = 0.9924990944342341
$ws.Range("K17").Value = 0.9923748021635791
$ws.Range("L17").Value = 0.9923570944068976
$ws.Range("M17").Value = 0.992346149021201

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "HexGrid-90degTilt22p5degRes"
$ws.Range("C18").Value = 0.9922571351621359
$ws.Range("D18").Value = 0.9940479545672992
$ws.Range("E18").Value = 0.9922574907241658
$ws.Range("F18").Value = 0.9922571351621359
$ws.Range("G18").Value = 0.9925161283939087
$ws.Range("H18").Value = 0.9900025626663325
$ws.Range("I18").Value = 0.9925469071083745
$ws.Range("J18").Value = 0.9940479545672992
$ws.Range("K18").Value = 0.9931527226457325
$ws.Range("L18").Value = 0.9927049289039342
$ws.Range("M18").Value = 0.9922713631037027

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C19").Value = 0.991251694699639
$ws.Range("D19").Value = 0.9972290896117837
$ws.Range("E19").Value = 0.9910574353791483
$ws.Range("F19").Value = 0.991251694699639
$ws.Range("G19").Value = 0.995299557398898
$ws.Range("H19").Value = 0.9882809289270912
$ws.Range("I19").Value = 0.9915358761318263
$ws.Range("J19").Value = 0.9972290896117837
$ws.Range("K19").Value = 0.994143262495466
$ws.Range("L19").Value = 0.9926974785975524
$ws.Range("M19").Value = 0.9924424303580643

# The three brand-new rows (17-19) need the same bold/boxed number style that
# column A uses for every other HKL index; copy it down from the row above.
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
